$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4607.4165
$ws.Range("I64").Value = 4473.625
$ws.Range("K64").Value = 4473.625
$ws.Range("M64").Value = -4225.625
$ws.Range("H67").Value = 4607.4165
$ws.Range("I67").Value = 4473.625
$ws.Range("K67").Value = 4473.625
$ws.Range("M67").Value = -3615.625
$ws.Range("H69").Value = 22272.727
$ws.Range("J69").Value = 22272.727
$ws.Range("L69").Value = 66818.181
$ws.Range("N69").Value = -68566.181
$ws.Range("H72").Value = 22272.727
$ws.Range("J72").Value = 22272.727
$ws.Range("L72").Value = 200454.543
$ws.Range("N72").Value = -209190.543
$ws.Range("H75").Value = 40314
$ws.Range("J75").Value = 40314
$ws.Range("L75").Value = 40314
$ws.Range("N75").Value = -42186
$ws.Range("H76").Value = 2621.25
$ws.Range("I76").Value = 2828.3333
$ws.Range("J76").Value = 2000
$ws.Range("K76").Value = 2828.3333
$ws.Range("L76").Value = 2000
$ws.Range("M76").Value = -2513.3333
$ws.Range("N76").Value = -2630
$ws.Range("H78").Value = 40314
$ws.Range("J78").Value = 40314
$ws.Range("L78").Value = 120942
$ws.Range("N78").Value = -130302
$ws.Range("H79").Value = 2621.25
$ws.Range("I79").Value = 2828.3333
$ws.Range("J79").Value = 2000
$ws.Range("K79").Value = 2828.3333
$ws.Range("L79").Value = 2000
$ws.Range("M79").Value = -1736.3333
$ws.Range("N79").Value = -4184
$ws.Range("H100").Value = 2477.2307
$ws.Range("I100").Value = 1900.4
$ws.Range("K100").Value = 1900.4
$ws.Range("M100").Value = -1359.4
$ws.Range("H132").Value = 4312
$ws.Range("I132").Value = 4312
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12936
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10406
$ws.Range("N132").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H55").Value = 29990
$ws.Range("J55").Value = 29990
$ws.Range("L55").Value = 29990
$ws.Range("N55").Value = -30620
$ws.Range("H132").Value = 7696290
$ws.Range("J132").Value = 9333
$ws.Range("L132").Value = 27999
$ws.Range("N132").Value = -33059

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 145.42857
$ws.Range("I11").Value = 32
$ws.Range("J11").Value = 296.66666
$ws.Range("K11").Value = 32
$ws.Range("L11").Value = 296.66666
$ws.Range("M11").Value = 108
$ws.Range("N11").Value = -576.66666
$ws.Range("H86").Value = 1880.0526
$ws.Range("I86").Value = 1947.875
$ws.Range("J86").Value = 1518.3334
$ws.Range("K86").Value = 1947.875
$ws.Range("L86").Value = 1518.3334
$ws.Range("M86").Value = -824.875
$ws.Range("N86").Value = -3764.3334
$ws.Range("H89").Value = 1880.0526
$ws.Range("I89").Value = 1947.875
$ws.Range("J89").Value = 1518.3334
$ws.Range("K89").Value = 9739.375
$ws.Range("L89").Value = 7591.666999999999
$ws.Range("M89").Value = -4123.375
$ws.Range("N89").Value = -18823.667
$ws.Range("H111").Value = 45999
$ws.Range("J111").Value = 45999
$ws.Range("L111").Value = 45999
$ws.Range("N111").Value = -54179
$ws.Range("H112").Value = 62999
$ws.Range("J112").Value = 62999
$ws.Range("L112").Value = 62999
$ws.Range("N112").Value = -65953
$ws.Range("H117").Value = 32989
$ws.Range("J117").Value = 32989
$ws.Range("L117").Value = 32989
$ws.Range("N117").Value = -42167
$ws.Range("H118").Value = 185999.5
$ws.Range("J118").Value = 185999.5
$ws.Range("L118").Value = 185999.5
$ws.Range("N118").Value = -189313.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1468.4
$ws.Range("I19").Value = 1435.5
$ws.Range("J19").Value = 1600
$ws.Range("K19").Value = 1435.5
$ws.Range("L19").Value = 1600
$ws.Range("M19").Value = -1265.5
$ws.Range("N19").Value = -1940
$ws.Range("H24").Value = 1468.4
$ws.Range("I24").Value = 1435.5
$ws.Range("J24").Value = 1600
$ws.Range("K24").Value = 1435.5
$ws.Range("L24").Value = 1600
$ws.Range("M24").Value = -1265.5
$ws.Range("N24").Value = -1940
$ws.Range("H74").Value = 45386.273
$ws.Range("I74").Value = 35000
$ws.Range("J74").Value = 47694.332
$ws.Range("K74").Value = 35000
$ws.Range("L74").Value = 47694.332
$ws.Range("M74").Value = -34126
$ws.Range("N74").Value = -49442.332
$ws.Range("H77").Value = 45386.273
$ws.Range("I77").Value = 35000
$ws.Range("J77").Value = 47694.332
$ws.Range("K77").Value = 105000
$ws.Range("L77").Value = 143082.996
$ws.Range("M77").Value = -100632
$ws.Range("N77").Value = -151818.996
$ws.Range("H86").Value = 13049.462
$ws.Range("J86").Value = 15549.875
$ws.Range("L86").Value = 15549.875
$ws.Range("N86").Value = -17795.875
$ws.Range("H89").Value = 13049.462
$ws.Range("J89").Value = 15549.875
$ws.Range("L89").Value = 77749.375
$ws.Range("N89").Value = -88981.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 263.75
$ws.Range("I6").Value = 185.33333
$ws.Range("J6").Value = 499
$ws.Range("K6").Value = 555.99999
$ws.Range("L6").Value = 1497
$ws.Range("M6").Value = -442.99999
$ws.Range("N6").Value = -1723
$ws.Range("H51").Value = 2515
$ws.Range("I51").Value = 2515
$ws.Range("K51").Value = 7545
$ws.Range("M51").Value = -7085
$ws.Range("H114").Value = 251742
$ws.Range("I114").Value = 335496
$ws.Range("J114").Value = 480
$ws.Range("K114").Value = 1006488
$ws.Range("L114").Value = 1440
$ws.Range("M114").Value = -1003234
$ws.Range("N114").Value = -7948
$ws.Range("H132").Value = 2103.5
$ws.Range("I132").Value = 1495
$ws.Range("J132").Value = 2171.111
$ws.Range("K132").Value = 13455
$ws.Range("L132").Value = 19539.999
$ws.Range("M132").Value = -10925
$ws.Range("N132").Value = -24599.999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 151.21428
$ws.Range("I2").Value = 101.75
$ws.Range("J2").Value = 217.16667
$ws.Range("K2").Value = 101.75
$ws.Range("L2").Value = 217.16667
$ws.Range("M2").Value = 11.25
$ws.Range("N2").Value = -443.16667
$ws.Range("H80").Value = 2314.6
$ws.Range("I80").Value = 2064
$ws.Range("J80").Value = 2422
$ws.Range("K80").Value = 2064
$ws.Range("L80").Value = 2422
$ws.Range("M80").Value = -1066
$ws.Range("N80").Value = -4418
$ws.Range("H83").Value = 2314.6
$ws.Range("I83").Value = 2064
$ws.Range("J83").Value = 2422
$ws.Range("K83").Value = 10320
$ws.Range("L83").Value = 12110
$ws.Range("M83").Value = -5328
$ws.Range("N83").Value = -22094
$ws.Range("H102").Value = 4376.533
$ws.Range("I102").Value = 3524.6155
$ws.Range("K102").Value = 3524.6155
$ws.Range("M102").Value = -1902.6155
$ws.Range("H107").Value = 4994.9165
$ws.Range("I107").Value = 3848.2222
$ws.Range("J107").Value = 8435
$ws.Range("K107").Value = 3848.2222
$ws.Range("L107").Value = 8435
$ws.Range("M107").Value = -1928.2222
$ws.Range("N107").Value = -12275

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 7558.6
$ws.Range("J30").Value = 17499.5
$ws.Range("L30").Value = 17499.5
$ws.Range("N30").Value = -17715.5
$ws.Range("H35").Value = 5983
$ws.Range("I35").Value = 1179.8
$ws.Range("K35").Value = 1179.8
$ws.Range("M35").Value = -843.8
$ws.Range("H55").Value = 586.73334
$ws.Range("I55").Value = 420.3
$ws.Range("K55").Value = 420.3
$ws.Range("M55").Value = -247.3
$ws.Range("H61").Value = 6700.625
$ws.Range("I61").Value = 6700.625
$ws.Range("K61").Value = 6700.625
$ws.Range("M61").Value = -6498.625
$ws.Range("H81").Value = 56998.5
$ws.Range("J81").Value = 56998.5
$ws.Range("L81").Value = 56998.5
$ws.Range("N81").Value = -58994.5
$ws.Range("H82").Value = 1347.8572
$ws.Range("I82").Value = 1450.1052
$ws.Range("K82").Value = 1450.1052
$ws.Range("M82").Value = -1089.1052
$ws.Range("H84").Value = 56998.5
$ws.Range("J84").Value = 56998.5
$ws.Range("L84").Value = 170995.5
$ws.Range("N84").Value = -180979.5
$ws.Range("H85").Value = 1347.8572
$ws.Range("I85").Value = 1450.1052
$ws.Range("K85").Value = 1450.1052
$ws.Range("M85").Value = -202.1052
$ws.Range("H93").Value = 2676.2727
$ws.Range("I93").Value = 1380.4
$ws.Range("K93").Value = 1380.4
$ws.Range("M93").Value = -132.4000000000001
$ws.Range("H100").Value = 12477344
$ws.Range("I100").Value = 14259286
$ws.Range("J100").Value = 3750
$ws.Range("K100").Value = 14259286
$ws.Range("L100").Value = 3750
$ws.Range("M100").Value = -14258745
$ws.Range("N100").Value = -4832
$ws.Range("H113").Value = 6700.625
$ws.Range("I113").Value = 6700.625
$ws.Range("K113").Value = 6700.625
$ws.Range("M113").Value = -4530.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 728.8
$ws.Range("I107").Value = 728.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2186.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -266.3999999999996
$ws.Range("N107").ClearContents()
